$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>37 x 13</w:t><w:br/><w:t xml:space="preserve">  1    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>'
$t.Cell(1,1).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>25 x 14</w:t><w:br/><w:t xml:space="preserve">  1    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>'
$t.Cell(1,2).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>43 x 20</w:t><w:br/><w:t xml:space="preserve">  2    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>'
$t.Cell(1,3).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>99 x 92</w:t><w:br/><w:t xml:space="preserve">  9    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p>'
$t.Cell(2,1).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>72 x 17</w:t><w:br/><w:t xml:space="preserve">  1    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p>'
$t.Cell(2,2).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>14 x 72</w:t><w:br/><w:t xml:space="preserve">  7    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p>'
$t.Cell(2,3).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>47 x 75</w:t><w:br/><w:t xml:space="preserve">  7    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>'
$t.Cell(3,1).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>33 x 14</w:t><w:br/><w:t xml:space="preserve">  1    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>'
$t.Cell(3,2).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>59 x 68</w:t><w:br/><w:t xml:space="preserve">  6    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p>'
$t.Cell(3,3).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>65 x 13</w:t><w:br/><w:t xml:space="preserve">  1    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>'
$t.Cell(4,1).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>12 x 43</w:t><w:br/><w:t xml:space="preserve">  4    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p>'
$t.Cell(4,2).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>94 x 87</w:t><w:br/><w:t xml:space="preserve">  8    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p>'
$t.Cell(4,3).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>10 x 76</w:t><w:br/><w:t xml:space="preserve">  7    6</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p>'
$t.Cell(5,1).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>44 x 14</w:t><w:br/><w:t xml:space="preserve">  1    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p>'
$t.Cell(5,2).Range.InsertXML($xml) | Out-Null

$xml = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>95 x 16</w:t><w:br/><w:t xml:space="preserve">  1    6</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>'
$t.Cell(5,3).Range.InsertXML($xml) | Out-Null
